$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: strikeouts column "K" (column G) now uses the
# actual K count instead of the old "Strike#" value. Update rows 2-11.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 3
$ws.Range("G10").Value = 9
$ws.Range("G11").Value = 1
